$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D whose new value looks like a plain number; Excel would
# otherwise auto-convert them from text to a numeric type. Force them to stay
# text (matching the original inlineStr text cells), then restore the default
# "Normal" style so no stray formatting is introduced.
$numericLookingDCells = @(
    'D5', 'D8', 'D11', 'D14', 'D19', 'D22', 'D25', 'D26', 'D28', 'D43', 'D44', 'D45', 'D48'
)
foreach ($ref in $numericLookingDCells) {
    $ws.Range($ref).NumberFormat = "@"
}

$ws.Range("D2").Value = '26.964.04'
$ws.Range("E2").Value = '  +0.11%  '

$ws.Range("D3").Value = '1.560.59'
$ws.Range("E3").Value = '  +0.44%  '

$ws.Range("D5").Value = '207.28'
$ws.Range("E5").Value = '  +0.12%  '

$ws.Range("E6").Value = '  +0.34%  '

$ws.Range("D8").Value = '22.10'
$ws.Range("E8").Value = '  +1.47%  '

$ws.Range("E9").Value = '  +0.21%  '

$ws.Range("E10").Value = '  +2.23%  '

$ws.Range("D11").Value = '0.0857'
$ws.Range("E11").Value = '  -0.08%  '

$ws.Range("D12").Value = '1.783.08'
$ws.Range("E12").Value = '  +0.44%  '

$ws.Range("D13").Value = '1.560.75'
$ws.Range("E13").Value = '  +0.40%  '

$ws.Range("D14").Value = '3.76'

$ws.Range("E15").Value = '  +0.40%  '

$ws.Range("E16").Value = '  +0.71%  '

$ws.Range("D17").Value = '26.973.97'
$ws.Range("E17").Value = '  +0.17%  '

$ws.Range("B18").Value = 'ShibaInu'
$ws.Range("C18").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D18").Value = '0.0₃0705'
$ws.Range("E18").Value = '  +2.36%  '

$ws.Range("B19").Value = 'BitcoinCash'
$ws.Range("C19").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D19").Value = '217.04'
$ws.Range("E19").Value = '  -0.01%  '

$ws.Range("E20").Value = '  +1.76%  '

$ws.Range("E21").Value = '  -0.23%  '

$ws.Range("D22").Value = '4.10'
$ws.Range("E22").Value = '  +1.58%  '

$ws.Range("E23").Value = '  -0.11%  '

$ws.Range("E24").Value = '  -1.43%  '

$ws.Range("D25").Value = '152.88'
$ws.Range("E25").Value = '  -0.68%  '

$ws.Range("D26").Value = '6.61'
$ws.Range("E26").Value = '  +0.66%  '

$ws.Range("E27").Value = '  +1.59%  '

$ws.Range("D28").Value = '0.104'
$ws.Range("E28").Value = '  +1.44%  '

$ws.Range("E29").Value = '  -0.31%  '

$ws.Range("E31").Value = '  +1.45%  '

$ws.Range("E32").Value = '  +0.70%  '

$ws.Range("E33").Value = '  +3.25%  '

$ws.Range("D34").Value = '1.420.70'
$ws.Range("E34").Value = '  -0.64%  '

$ws.Range("E35").Value = '  +3.34%  '

$ws.Range("E36").Value = '  +9.35%  '

$ws.Range("E37").Value = '  +1.51%  '

$ws.Range("E38").Value = '  +0.69%  '

$ws.Range("E39").Value = '  +2.36%  '

$ws.Range("E40").Value = '  -0.26%  '

$ws.Range("E41").Value = '  -0.24%  '

$ws.Range("E42").Value = '  +0.13%  '

$ws.Range("B43").Value = 'MXToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D43").Value = '2.33'
$ws.Range("E43").Value = '  +2.31%  '

$ws.Range("B44").Value = 'WEMIXToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D44").Value = '1.01'
$ws.Range("E44").Value = '  +1.84%  '

$ws.Range("D45").Value = '64.72'
$ws.Range("E45").Value = '  +1.47%  '

$ws.Range("E46").Value = '  +0.55%  '

$ws.Range("D47").Value = '1.696.65'
$ws.Range("E47").Value = '  +0.40%  '

$ws.Range("D48").Value = '87.57'
$ws.Range("E48").Value = '  +1.61%  '

$ws.Range("E49").Value = '  -0.49%  '

$ws.Range("D50").Value = '0.0₆0100'
$ws.Range("E50").Value = '  +0.20%  '

$ws.Range("E51").Value = '  +0.15%  '

foreach ($ref in $numericLookingDCells) {
    $ws.Range($ref).Style = "Normal"
}
